# Update "苏州-漫展信息" workbook to the next scrape snapshot.
#
# The underlying event list dropped the "苏州·无限次元夜场" entry (it fell off
# the upstream feed) and several "想去人数" (interest count) figures ticked up
# for events still listed. This touches the two sheets that contained that
# row: "展览" (Exhibitions) and "全部类型" (All types). "演出" and "本地生活"
# are untouched.

$wb = $excel.ActiveWorkbook

# Sheet name -> list of (Link URL, new F-column "want to go" count) updates
# that apply to rows surviving the deletion. Cross-references the row by its
# unique Link (column H) so it doesn't matter that the row number shifts up
# by one once the "无限次元夜场" row is removed.
$fUpdates = @{
    "https://show.bilibili.com/platform/detail.html?id=91280" = 818
    "https://show.bilibili.com/platform/detail.html?id=93201" = 1172
    "https://show.bilibili.com/platform/detail.html?id=93516" = 16
    "https://show.bilibili.com/platform/detail.html?id=91626" = 287
    "https://show.bilibili.com/platform/detail.html?id=93335" = 1030
    "https://show.bilibili.com/platform/detail.html?id=91709" = 526
    "https://show.bilibili.com/platform/detail.html?id=91324" = 553
    "https://show.bilibili.com/platform/detail.html?id=92177" = 13157
    "https://show.bilibili.com/platform/detail.html?id=93234" = 5401
    "https://show.bilibili.com/platform/detail.html?id=84858" = 5556
    "https://show.bilibili.com/platform/detail.html?id=93262" = 14
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 5 is "苏州·无限次元夜场" (2024-10-19) in both sheets - find it
    # defensively by scanning column C rather than assuming the row index,
    # then delete the whole row; Excel shifts everything below it up.
    $lastRow = $ws.UsedRange.Rows.Count
    $targetRow = 0
    for ($r = 2; $r -le $lastRow; $r++) {
        $name = $ws.Cells.Item($r, 3).Value()
        if ($name -eq "苏州·无限次元夜场") {
            $targetRow = $r
            break
        }
    }
    if ($targetRow -gt 0) {
        $ws.Rows.Item($targetRow).Delete()
    }

    # Renumber column A (the literal running index 1..N) now that a row is
    # gone, and refresh the "want to go" counts that moved on.
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1

        $link = $ws.Cells.Item($r, 8).Value()
        if ($fUpdates.ContainsKey($link)) {
            $ws.Cells.Item($r, 6).Value = $fUpdates[$link]
        }
    }
}
